$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) column C for rows 2-28 from 45554 -> 45555
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45555
}

# Ensure row 28 carries an explicit row height (matches the diff: ht="15" customHeight="1")
$ws.Rows.Item(28).RowHeight = 15

# Add the new data row 29
$ws.Cells.Item(29, 1).Value = "A 38636-2024"
$ws.Cells.Item(29, 2).Value = 45546
$ws.Cells.Item(29, 3).Value = 45555
$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"
$ws.Cells.Item(29, 7).Value = 0.7
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0

# Apply the same date number format as B28/C28 to B29/C29
$ws.Cells.Item(29, 2).NumberFormat = $ws.Cells.Item(28, 2).NumberFormat
$ws.Cells.Item(29, 3).NumberFormat = $ws.Cells.Item(28, 3).NumberFormat

# R29 mirrors R28: an empty, wrap-text-styled cell
$ws.Cells.Item(29, 18).WrapText = $true
